$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with freshly scraped crypto data.
# D-column values are forced to Text format so strings like "604.42" or
# "64.478.43" are preserved exactly (not reinterpreted as numbers/dates).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.478.43"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.51%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.141.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.69%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "604.42"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.29%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.33"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.53%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.136.26"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.79%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.15%  "

$ws.Range("E10").Value = "  +1.60%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.40"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.25%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.43%  "

$ws.Range("E13").Value = "  +4.35%  "

$ws.Range("E14").Value = "  +0.98%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.658.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.70%  "

$ws.Range("E16").Value = "  +3.05%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.389.69"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.173.82"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.95%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.88"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.69%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "478.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.48%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.66%  "

$ws.Range("E22").Value = "  +2.07%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.69"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.87%  "

$ws.Range("E24").Value = "  +2.38%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.43"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.03%  "

$ws.Range("E26").Value = "  -0.11%  "

$ws.Range("E27").Value = "  -0.27%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.40"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.28%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.25"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.15%  "

$ws.Range("E30").Value = "  -4.01%  "

$ws.Range("E31").Value = "  +3.39%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "27.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +4.04%  "

$ws.Range("E33").Value = "  +0.06%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.66"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.21%  "

$ws.Range("E35").Value = "  +0.72%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0₃0778"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +7.92%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.98"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.82%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.97%  "

$ws.Range("E39").Value = "  +5.14%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "444.76"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -2.40%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0394"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.05%  "

$ws.Range("E42").Value = "  +1.03%  "

$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.851.16"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.48%  "

$ws.Range("E45").Value = "  -0.47%  "

$ws.Range("E46").Value = "  +0.58%  "

$ws.Range("E47").Value = "  +2.28%  "

$ws.Range("E48").Value = "  +0.05%  "

$ws.Range("E49").Value = "  +0.75%  "

$ws.Range("E50").Value = "  +0.78%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.10"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.07%  "
